$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 99
$ws.Range("H99").Value = 1863.4286
$ws.Range("I99").Value = 1884
$ws.Range("K99").Value = 5652
$ws.Range("M99").Value = -4154
# Row 107
$ws.Range("H107").Value = 359.22726
$ws.Range("I107").Value = 328.7143
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 328.7143
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1591.2857
$ws.Range("N107").Value = -4840
# Row 111
$ws.Range("H111").Value = 1716.9412
$ws.Range("I111").Value = 1439.75
$ws.Range("J111").Value = 2382.2
$ws.Range("K111").Value = 4319.25
$ws.Range("L111").Value = 7146.599999999999
$ws.Range("M111").Value = -1252.25
$ws.Range("N111").Value = -13280.6

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 77.666664
$ws.Range("I5").Value = 79.61539
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 79.61539
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = 32.38461
$ws.Range("N5").Value = -289
# Row 88
$ws.Range("H88").Value = 2693.077
$ws.Range("I88").Value = 2550
$ws.Range("J88").Value = 2782.5
$ws.Range("K88").Value = 2550
$ws.Range("L88").Value = 2782.5
$ws.Range("M88").Value = -2144
$ws.Range("N88").Value = -3594.5
# Row 91
$ws.Range("H91").Value = 2693.077
$ws.Range("I91").Value = 2550
$ws.Range("J91").Value = 2782.5
$ws.Range("K91").Value = 2550
$ws.Range("L91").Value = 2782.5
$ws.Range("M91").Value = -1146
$ws.Range("N91").Value = -5590.5

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 77.666664
$ws.Range("I4").Value = 79.61539
$ws.Range("J4").Value = 65
$ws.Range("K4").Value = 79.61539
$ws.Range("L4").Value = 65
$ws.Range("M4").Value = 35.38461
$ws.Range("N4").Value = -295
# Row 86
$ws.Range("H86").Value = 373569.53
$ws.Range("I86").Value = 2305.2666
$ws.Range("J86").Value = 837649.8
$ws.Range("K86").Value = 2305.2666
$ws.Range("L86").Value = 837649.8
$ws.Range("M86").Value = -1182.2666
$ws.Range("N86").Value = -839895.8
# Row 89
$ws.Range("H89").Value = 373569.53
$ws.Range("I89").Value = 2305.2666
$ws.Range("J89").Value = 837649.8
$ws.Range("K89").Value = 11526.333
$ws.Range("L89").Value = 4188249
$ws.Range("M89").Value = -5910.332999999999
$ws.Range("N89").Value = -4199481
# Row 99
$ws.Range("H99").Value = 991.3333
$ws.Range("I99").Value = 873.7857
$ws.Range("J99").Value = 1402.75
$ws.Range("K99").Value = 873.7857
$ws.Range("L99").Value = 1402.75
$ws.Range("M99").Value = 624.2143
$ws.Range("N99").Value = -4398.75

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2770.375
$ws.Range("I16").Value = 1314.8182
$ws.Range("J16").Value = 5972.6
$ws.Range("K16").Value = 1314.8182
$ws.Range("L16").Value = 5972.6
$ws.Range("M16").Value = -1027.8182
$ws.Range("N16").Value = -6546.6
# Row 31
$ws.Range("H31").Value = 23450.076
$ws.Range("I31").Value = 1243.0322
$ws.Range("J31").Value = 56231.906
$ws.Range("K31").Value = 1243.0322
$ws.Range("L31").Value = 56231.906
$ws.Range("M31").Value = -948.0322000000001
$ws.Range("N31").Value = -56821.906
# Row 34
$ws.Range("H34").Value = 23450.076
$ws.Range("I34").Value = 1243.0322
$ws.Range("J34").Value = 56231.906
$ws.Range("K34").Value = 1243.0322
$ws.Range("L34").Value = 56231.906
$ws.Range("M34").Value = -1041.0322
$ws.Range("N34").Value = -56635.906
# Row 62
$ws.Range("H62").Value = 4496.393
$ws.Range("I62").Value = 4882.609
$ws.Range("J62").Value = 2719.8
$ws.Range("K62").Value = 4882.609
$ws.Range("L62").Value = 2719.8
$ws.Range("M62").Value = -4258.609
$ws.Range("N62").Value = -3967.8
# Row 65
$ws.Range("H65").Value = 4496.393
$ws.Range("I65").Value = 4882.609
$ws.Range("J65").Value = 2719.8
$ws.Range("K65").Value = 24413.045
$ws.Range("L65").Value = 13599
$ws.Range("M65").Value = -21293.045
$ws.Range("N65").Value = -19839
# Row 99
$ws.Range("H99").Value = 1347.7587
$ws.Range("I99").Value = 1314.55
$ws.Range("J99").Value = 1421.5555
$ws.Range("K99").Value = 1314.55
$ws.Range("L99").Value = 1421.5555
$ws.Range("M99").Value = 183.45
$ws.Range("N99").Value = -4417.5555
# Row 107
$ws.Range("H107").Value = 1018
$ws.Range("I107").Value = 1039.8
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1039.8
$ws.Range("L107").Value = 800
$ws.Range("M107").Value = 880.2
$ws.Range("N107").Value = -4640
# Row 113
$ws.Range("H113").Value = 2770.375
$ws.Range("I113").Value = 1314.8182
$ws.Range("J113").Value = 5972.6
$ws.Range("K113").Value = 1314.8182
$ws.Range("L113").Value = 5972.6
$ws.Range("M113").Value = 855.1818000000001
$ws.Range("N113").Value = -10312.6
# Row 126
$ws.Range("H126").Value = 1347.7587
$ws.Range("I126").Value = 1314.55
$ws.Range("J126").Value = 1421.5555
$ws.Range("K126").Value = 3943.65
$ws.Range("L126").Value = 4264.666499999999
$ws.Range("M126").Value = -1473.65
$ws.Range("N126").Value = -9204.666499999999
# Row 132
$ws.Range("H132").Value = 50007620
$ws.Range("I132").Value = 83343990
$ws.Range("J132").Value = 3062.75
$ws.Range("K132").Value = 250031970
$ws.Range("L132").Value = 9188.25
$ws.Range("M132").Value = -250029440
$ws.Range("N132").Value = -14248.25

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 600.6667
$ws.Range("I5").Value = 423.65216
$ws.Range("J5").Value = 913.8461
$ws.Range("K5").Value = 1270.95648
$ws.Range("L5").Value = 2741.5383
$ws.Range("M5").Value = -1158.95648
$ws.Range("N5").Value = -2965.5383
# Row 64
$ws.Range("H64").Value = 83335816
$ws.Range("I64").Value = 3939
$ws.Range("J64").Value = 111113100
$ws.Range("K64").Value = 11817
$ws.Range("L64").Value = 333339300
$ws.Range("M64").Value = -11547
$ws.Range("N64").Value = -333339840
# Row 67
$ws.Range("H67").Value = 83335816
$ws.Range("I67").Value = 3939
$ws.Range("J67").Value = 111113100
$ws.Range("K67").Value = 11817
$ws.Range("L67").Value = 333339300
$ws.Range("M67").Value = -10881
$ws.Range("N67").Value = -333341172
# Row 122
$ws.Range("H122").Value = 1030.25
$ws.Range("J122").Value = 1581
$ws.Range("L122").Value = 14229
$ws.Range("N122").Value = -19129
# Row 132
$ws.Range("H132").Value = 1063.9678
$ws.Range("I132").Value = 892.1111
$ws.Range("K132").Value = 8028.9999
$ws.Range("M132").Value = -5498.9999
# Row 135
$ws.Range("H135").Value = 600.6667
$ws.Range("I135").Value = 423.65216
$ws.Range("J135").Value = 913.8461
$ws.Range("K135").Value = 3812.86944
$ws.Range("L135").Value = 8224.6149
$ws.Range("M135").Value = -1277.86944
$ws.Range("N135").Value = -13294.6149

$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 21838774
$ws.Range("I11").Value = 40002500
$ws.Range("J11").Value = 42302.4
$ws.Range("K11").Value = 40002500
$ws.Range("L11").Value = 42302.4
$ws.Range("N11").Value = -42580.4
$ws.Range("M11").Value = -40002361
# Row 51
$ws.Range("H51").Value = 36000
$ws.Range("J51").Value = 36000
$ws.Range("L51").Value = 36000
$ws.Range("N51").Value = -37018
# Row 70
$ws.Range("H70").Value = 4827.8667
$ws.Range("I70").Value = 5037.5
$ws.Range("J70").Value = 4588.2856
$ws.Range("K70").Value = 5037.5
$ws.Range("L70").Value = 4588.2856
$ws.Range("M70").Value = -4767.5
$ws.Range("N70").Value = -5128.2856
# Row 73
$ws.Range("H73").Value = 4827.8667
$ws.Range("I73").Value = 5037.5
$ws.Range("J73").Value = 4588.2856
$ws.Range("K73").Value = 5037.5
$ws.Range("L73").Value = 5037.5
$ws.Range("M73").Value = -4101.5
$ws.Range("N73").Value = -6460.2856
# Row 113
$ws.Range("H113").Value = 3292.2144
$ws.Range("I113").Value = 3066.7778
$ws.Range("J113").Value = 3698
$ws.Range("K113").Value = 3066.7778
$ws.Range("L113").Value = 3698
$ws.Range("M113").Value = -896.7777999999998
$ws.Range("N113").Value = -8038
# Row 126
$ws.Range("H126").Value = 1900
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1900
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5700
$ws.Range("N126").Value = -10640
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 6
$ws.Range("H6").Value = 27980
$ws.Range("J6").Value = 27980
$ws.Range("L6").Value = 27980
$ws.Range("N6").Value = -28204
# Row 40
$ws.Range("H40").Value = 4231.1875
$ws.Range("I40").Value = 3955.5557
$ws.Range("J40").Value = 4585.5713
$ws.Range("K40").Value = 3955.5557
$ws.Range("L40").Value = 4585.5713
$ws.Range("M40").Value = -3819.5557
$ws.Range("N40").Value = -4857.5713
